$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.143.81"
$ws.Range("E2").Value = "  -2.93%  "
$ws.Range("D3").Value = "3.013.86"
$ws.Range("E3").Value = "  -5.17%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'564.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.31%  "
$ws.Range("D6").Value = "'128.64"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.96%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.017.95"
$ws.Range("E8").Value = "  -4.99%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").Value = "'0.135"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.53%  "
$ws.Range("D11").Value = "'5.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "'0.434"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.41%  "
$ws.Range("D13").Value = "'0.0000224"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.06%  "
$ws.Range("D14").Value = "'33.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.33%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").Value = "3.525.02"
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").Value = "61.397.05"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "3.016.14"
$ws.Range("E18").Value = "  -4.86%  "
$ws.Range("D19").Value = "'6.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.72%  "
$ws.Range("D20").Value = "'441.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.19%  "
$ws.Range("D21").Value = "'13.20"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.53%  "
$ws.Range("D22").Value = "'0.666"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.73%  "
$ws.Range("D23").Value = "'7.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.44%  "
$ws.Range("D24").Value = "'12.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.90%  "
$ws.Range("D25").Value = "'79.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.83%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'2.50"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.21%  "
$ws.Range("D29").Value = "'7.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.41%  "
$ws.Range("D30").Value = "'1.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.86%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'25.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.54%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.97%  "
$ws.Range("D33").Value = "'0.0941"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.48%  "
$ws.Range("D34").Value = "'2.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D35").Value = "'0.959"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.77%  "
$ws.Range("D36").Value = "'5.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.24%  "
$ws.Range("D37").Value = "'50.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").Value = "0.0₃0675"
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("D39").Value = "'0.0360"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.08%  "
$ws.Range("D40").Value = "'7.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'380.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.03%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.108"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.69%  "
$ws.Range("D43").Value = "2.688.11"
$ws.Range("E43").Value = "  -3.83%  "
$ws.Range("D44").Value = "'2.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.27%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.236"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.14%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'34.05"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").Value = "'120.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").Value = "'1.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.04%  "
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").Value = "'23.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.14%  "
